# Update market/profit data on several sheets (Gungnir scheduled runner refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 78009.234
$ws.Range("I17").Value = 800
$ws.Range("J17").Value = 92047.27
$ws.Range("K17").Value = 2400
$ws.Range("L17").Value = 276141.81
$ws.Range("M17").Value = -2232
$ws.Range("N17").Value = -276477.81

$ws.Range("H93").Value = 37933.332
$ws.Range("J93").Value = 37933.332
$ws.Range("L93").Value = 37933.332
$ws.Range("N93").Value = -42925.332

$ws.Range("H113").Value = 3400
$ws.Range("I113").Value = 3750
$ws.Range("J113").Value = 3272.7273
$ws.Range("K113").Value = 3750
$ws.Range("L113").Value = 3272.7273
$ws.Range("M113").Value = -496
$ws.Range("N113").Value = -9780.7273

$ws.Range("H116").Value = 639.625
$ws.Range("I116").Value = 639.625
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 639.625
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 2802.375
$ws.Range("N116").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24401516
$ws.Range("I32").Value = 9714.839
$ws.Range("J32").Value = 100016100
$ws.Range("K32").Value = 9714.839
$ws.Range("L32").Value = 100016100
$ws.Range("M32").Value = -9427.839
$ws.Range("N32").Value = -100016674

$ws.Range("H63").Value = 2571.2856
$ws.Range("J63").Value = 3166.6667
$ws.Range("L63").Value = 3166.6667
$ws.Range("N63").Value = -4538.6667

$ws.Range("H66").Value = 2571.2856
$ws.Range("J66").Value = 3166.6667
$ws.Range("L66").Value = 15833.3335
$ws.Range("N66").Value = -22697.3335

$ws.Range("H122").Value = 1418.9565
$ws.Range("I122").Value = 1218.1538
$ws.Range("J122").Value = 1680
$ws.Range("K122").Value = 3654.4614
$ws.Range("L122").Value = 5040
$ws.Range("M122").Value = -1204.4614
$ws.Range("N122").Value = -9940

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 9715.357
$ws.Range("I82").Value = 2768.9
$ws.Range("J82").Value = 27081.5
$ws.Range("K82").Value = 2768.9
$ws.Range("L82").Value = 27081.5
$ws.Range("M82").Value = -2385.9
$ws.Range("N82").Value = -27847.5

$ws.Range("H85").Value = 9715.357
$ws.Range("I85").Value = 2768.9
$ws.Range("J85").Value = 27081.5
$ws.Range("K85").Value = 2768.9
$ws.Range("L85").Value = 27081.5
$ws.Range("M85").Value = -1442.9
$ws.Range("N85").Value = -29733.5

$ws.Range("H92").Value = 40000
$ws.Range("J92").Value = 40000
$ws.Range("L92").Value = 40000
$ws.Range("N92").Value = -44992

$ws.Range("H134").Value = 5563946
$ws.Range("I134").Value = 1778
$ws.Range("J134").Value = 6954488
$ws.Range("K134").Value = 5334
$ws.Range("L134").Value = 20863464
$ws.Range("M134").Value = -2799
$ws.Range("N134").Value = -20868534

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 43480010
$ws.Range("I58").Value = 142858270
$ws.Range("J58").Value = 2018.75
$ws.Range("K58").Value = 142858270
$ws.Range("L58").Value = 2018.75
$ws.Range("M58").Value = -142858067
$ws.Range("N58").Value = -2424.75

$ws.Range("H88").Value = 27414.334
$ws.Range("J88").Value = 27414.334
$ws.Range("L88").Value = 27414.334
$ws.Range("N88").Value = -28226.334

$ws.Range("H91").Value = 27414.334
$ws.Range("J91").Value = 27414.334
$ws.Range("L91").Value = 27414.334
$ws.Range("N91").Value = -30222.334

$ws.Range("H99").Value = 45456544
$ws.Range("I99").Value = 125001600
$ws.Range("J99").Value = 2224.2856
$ws.Range("K99").Value = 125001600
$ws.Range("L99").Value = 2224.2856
$ws.Range("M99").Value = -125000102
$ws.Range("N99").Value = -5220.2856

$ws.Range("H126").Value = 45456544
$ws.Range("I126").Value = 125001600
$ws.Range("J126").Value = 2224.2856
$ws.Range("K126").Value = 375004800
$ws.Range("L126").Value = 6672.8568
$ws.Range("M126").Value = -375002330
$ws.Range("N126").Value = -11612.8568

$ws.Range("H134").Value = 1222.5454
$ws.Range("I134").Value = 1183.1428
$ws.Range("J134").Value = 1291.5
$ws.Range("K134").Value = 3549.4284
$ws.Range("L134").Value = 3874.5
$ws.Range("M134").Value = -1014.4284
$ws.Range("N134").Value = -8944.5

$ws.Range("H136").Value = 43480010
$ws.Range("I136").Value = 142858270
$ws.Range("J136").Value = 2018.75
$ws.Range("K136").Value = 428574810
$ws.Range("L136").Value = 6056.25
$ws.Range("M136").Value = -428572260
$ws.Range("N136").Value = -11156.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 747.1900000000001
$ws.Range("J131").Value = 773.7527
$ws.Range("L131").Value = 2321.2581
$ws.Range("N131").Value = -12401.2581

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").Value = ""

$ws.Range("H122").Value = 50012820
$ws.Range("I122").Value = 83351030
$ws.Range("J122").Value = 5500.25
$ws.Range("K122").Value = 250053090
$ws.Range("L122").Value = 16500.75
$ws.Range("M122").Value = -250050640
$ws.Range("N122").Value = -21400.75

$ws.Range("H126").Value = 1770.5834
$ws.Range("I126").Value = 1283
$ws.Range("J126").Value = 3233.3333
$ws.Range("K126").Value = 3849
$ws.Range("L126").Value = 9699.999899999999
$ws.Range("M126").Value = -1379
$ws.Range("N126").Value = -14639.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 942.0417
$ws.Range("I7").Value = 965.3913
$ws.Range("J7").Value = 405
$ws.Range("K7").Value = 965.3913
$ws.Range("L7").Value = 405
$ws.Range("M7").Value = -853.3913
$ws.Range("N7").Value = -629

$ws.Range("H94").Value = 24330
$ws.Range("J94").Value = 24330
$ws.Range("L94").Value = 24330
$ws.Range("N94").Value = -25682

$ws.Range("H126").Value = 942.0417
$ws.Range("I126").Value = 965.3913
$ws.Range("J126").Value = 405
$ws.Range("K126").Value = 2896.1739
$ws.Range("L126").Value = 1215
$ws.Range("M126").Value = -426.1738999999998
$ws.Range("N126").Value = -6155

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 45660
$ws.Range("J46").Value = 45660
$ws.Range("L46").Value = 45660
$ws.Range("N46").Value = -46122

$ws.Range("H134").Value = 45660
$ws.Range("J134").Value = 45660
$ws.Range("L134").Value = 136980
$ws.Range("N134").Value = -142050
